$wb = $excel.ActiveWorkbook

# --- Update timestamps on the "data" sheet (F2:F5) ---
$dataSheet = $wb.Worksheets.Item("data")
$dataSheet.Range("F2").Value = "2021-10-05 14:19:03.173932"
$dataSheet.Range("F3").Value = "2021-10-05 14:19:03.173940"
$dataSheet.Range("F4").Value = "2021-10-05 14:19:03.173943"
$dataSheet.Range("F5").Value = "2021-10-05 14:19:03.173946"

# --- Add the new "metadata" sheet right after "data" ---
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# Header row
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Header styling (bold, centered, bordered) to mirror "data" sheet header style
$headerRange = $ws.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data row
$ws.Range("A2").Value = 0
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4160
$ws.Range("A2").Borders.LineStyle = 1

$ws.Range("B2").Value = "Additional findings health related - CNV analysis children"
$ws.Range("C2").Value = 932
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1.0"
$ws.Range("E2").Value = "2021-04-07T10:24:28.719980Z"
$ws.Range("F2").Value = "2021-10-05 14:19:03.170039"
$ws.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/932/?format=json"

$dataSheet.Select()
